$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# Shape 1: Title Placeholder 1 -> "Click to edit Master title style"
$titleShape = $m.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "{g0}ickclay otay edithay astermay itletay estylay{/g1}"

# Shape 2: Text Placeholder 2 -> 5 paragraphs (body levels 0-4)
$bodyShape = $m.Shapes.Item(2)
$bodyTextRange = $bodyShape.TextFrame.TextRange
$bodyTextRange.Paragraphs(1, 1).Text = "{g0}ickclay otay edithay astermay exttay esstylay{/g1}"
$bodyTextRange.Paragraphs(2, 1).Text = "{g0}econdsay evellay{/g1}"
$bodyTextRange.Paragraphs(3, 1).Text = "{g0}irdthay evellay{/g1}"
$bodyTextRange.Paragraphs(4, 1).Text = "{g0}ourthfay evellay{/g1}"
$bodyTextRange.Paragraphs(5, 1).Text = "{g0}ifthfay evellay{/g1}"

# Shape 3: Date Placeholder 3 -> "3/1/2007"
$dateShape = $m.Shapes.Item(3)
$dateShape.TextFrame.TextRange.Text = "{g0}3/1/2007{/g1}"

# Shape 5: Slide Number Placeholder 5 -> "‹#›"
$slideNumShape = $m.Shapes.Item(5)
$slideNumShape.TextFrame.TextRange.Text = "{g0}‹#›{/g1}"
